# Auto-generated edit script applying the Shiva_Profits.xlsx leve-profit refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 584.26666
$ws.Range("I53").Value = 751.125
$ws.Range("K53").Value = 751.125
$ws.Range("M53").Value = -114.125
$ws.Range("H107").Value = 18519346
$ws.Range("I107").Value = 25000748
$ws.Range("K107").Value = 25000748
$ws.Range("M107").Value = -24998828
$ws.Range("H112").Value = 2709.4412
$ws.Range("J112").Value = 2886.2222
$ws.Range("L112").Value = 8658.6666
$ws.Range("N112").Value = -10874.6666
$ws.Range("H113").Value = 3521.5652
$ws.Range("I113").Value = 2896.5833
$ws.Range("K113").Value = 2896.5833
$ws.Range("M113").Value = 357.4167000000002
$ws.Range("H116").Value = 5652.1665
$ws.Range("I116").Value = 6134.9
$ws.Range("J116").Value = 5048.75
$ws.Range("K116").Value = 6134.9
$ws.Range("L116").Value = 5048.75
$ws.Range("M116").Value = -2692.9
$ws.Range("N116").Value = -11932.75
$ws.Range("H129").Value = 1663.4
$ws.Range("J129").Value = 2055.3333
$ws.Range("L129").Value = 6165.999899999999
$ws.Range("N129").Value = -16165.9999
$ws.Range("H138").Value = 2795.859
$ws.Range("I138").Value = 2405.2
$ws.Range("K138").Value = 7215.599999999999
$ws.Range("M138").Value = -2075.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 3666.2222
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H74").Value = 3344.1
$ws.Range("J74").Value = 3917.1428
$ws.Range("L74").Value = 3917.1428
$ws.Range("N74").Value = -5665.1428
$ws.Range("H77").Value = 3344.1
$ws.Range("J77").Value = 3917.1428
$ws.Range("L77").Value = 19585.714
$ws.Range("N77").Value = -28321.714
$ws.Range("H102").Value = 1795.091
$ws.Range("I102").Value = 1455.8572
$ws.Range("J102").Value = 2388.75
$ws.Range("K102").Value = 1455.8572
$ws.Range("L102").Value = 2388.75
$ws.Range("M102").Value = 166.1428000000001
$ws.Range("N102").Value = -5632.75
$ws.Range("H125").Value = 115600
$ws.Range("J125").Value = 115600
$ws.Range("L125").Value = 115600
$ws.Range("N125").Value = -125440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2543.7778
$ws.Range("I105").Value = 2549.25
$ws.Range("K105").Value = 2549.25
$ws.Range("M105").Value = -802.25
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H134").Value = 2202.9614
$ws.Range("I134").Value = 1608.1628
$ws.Range("J134").Value = 5044.778
$ws.Range("K134").Value = 4824.4884
$ws.Range("L134").Value = 15134.334
$ws.Range("M134").Value = -2289.4884
$ws.Range("N134").Value = -20204.334
$ws.Range("H140").Value = 55366.11
$ws.Range("J140").Value = 55366.11
$ws.Range("L140").Value = 55366.11
$ws.Range("N140").Value = -65726.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 29546.5
$ws.Range("I22").Value = 29546.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 29546.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -29196.5
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 4984.6743
$ws.Range("I31").Value = 2544.1538
$ws.Range("J31").Value = 8717.235000000001
$ws.Range("K31").Value = 2544.1538
$ws.Range("L31").Value = 8717.235000000001
$ws.Range("M31").Value = -2249.1538
$ws.Range("N31").Value = -9307.235000000001
$ws.Range("H34").Value = 4984.6743
$ws.Range("I34").Value = 2544.1538
$ws.Range("J34").Value = 8717.235000000001
$ws.Range("K34").Value = 2544.1538
$ws.Range("L34").Value = 8717.235000000001
$ws.Range("M34").Value = -2342.1538
$ws.Range("N34").Value = -9121.235000000001
$ws.Range("H109").Value = 105047.164
$ws.Range("J109").Value = 105047.164
$ws.Range("L109").Value = 105047.164
$ws.Range("N109").Value = -107127.164
$ws.Range("H122").Value = 5397.4546
$ws.Range("I122").Value = 5876.55
$ws.Range("J122").Value = 606.5
$ws.Range("K122").Value = 17629.65
$ws.Range("L122").Value = 1819.5
$ws.Range("M122").Value = -15179.65
$ws.Range("N122").Value = -6719.5
$ws.Range("H134").Value = 3804.04
$ws.Range("I134").Value = 3500.3416
$ws.Range("J134").Value = 5187.5557
$ws.Range("K134").Value = 10501.0248
$ws.Range("L134").Value = 15562.6671
$ws.Range("M134").Value = -7966.024800000001
$ws.Range("N134").Value = -20632.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2030.1111
$ws.Range("I136").Value = 2030.1111
$ws.Range("K136").Value = 6090.3333
$ws.Range("M136").Value = -990.3333000000002
$ws.Range("H137").Value = 5166.3335
$ws.Range("I137").Value = 6499.3335
$ws.Range("J137").Value = 3833.3333
$ws.Range("K137").Value = 19498.0005
$ws.Range("L137").Value = 11499.9999
$ws.Range("M137").Value = -14398.0005
$ws.Range("N137").Value = -21699.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 21822432
$ws.Range("I80").Value = 40003050
$ws.Range("K80").Value = 40003050
$ws.Range("M80").Value = -40002052
$ws.Range("H83").Value = 21822432
$ws.Range("I83").Value = 40003050
$ws.Range("K83").Value = 200015250
$ws.Range("M83").Value = -200010258
$ws.Range("H97").Value = 1074.2727
$ws.Range("I97").Value = 1059.375
$ws.Range("J97").Value = 1114
$ws.Range("K97").Value = 1059.375
$ws.Range("L97").Value = 1114
$ws.Range("M97").Value = -563.375
$ws.Range("N97").Value = -2106
$ws.Range("H102").Value = 3778.6843
$ws.Range("I102").Value = 3646.3257
$ws.Range("J102").Value = 4185.2144
$ws.Range("K102").Value = 3646.3257
$ws.Range("L102").Value = 4185.2144
$ws.Range("M102").Value = -2024.3257
$ws.Range("N102").Value = -7429.2144
$ws.Range("H122").Value = 590108.5
$ws.Range("I122").Value = 1000973.6
$ws.Range("J122").Value = 3158.4285
$ws.Range("K122").Value = 3002920.8
$ws.Range("L122").Value = 9475.2855
$ws.Range("M122").Value = -3000470.8
$ws.Range("N122").Value = -14375.2855
$ws.Range("H123").Value = 26038.428
$ws.Range("J123").Value = 26038.428
$ws.Range("L123").Value = 26038.428
$ws.Range("N123").Value = -30938.428
$ws.Range("H126").Value = 4668.8237
$ws.Range("I126").Value = 3597.5557
$ws.Range("J126").Value = 5874
$ws.Range("K126").Value = 10792.6671
$ws.Range("L126").Value = 17622
$ws.Range("M126").Value = -8322.667099999999
$ws.Range("N126").Value = -22562
$ws.Range("H132").Value = 1843.4054
$ws.Range("I132").Value = 1174.1
$ws.Range("K132").Value = 3522.3
$ws.Range("M132").Value = -992.2999999999997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 12501412
$ws.Range("I16").Value = 20834572
$ws.Range("J16").Value = 1671
$ws.Range("K16").Value = 20834572
$ws.Range("L16").Value = 1671
$ws.Range("M16").Value = -20834402
$ws.Range("N16").Value = -2011
$ws.Range("H22").Value = 995.6
$ws.Range("J22").Value = 996.6667
$ws.Range("L22").Value = 996.6667
$ws.Range("N22").Value = -1586.6667
$ws.Range("H27").Value = 995.6
$ws.Range("J27").Value = 996.6667
$ws.Range("L27").Value = 996.6667
$ws.Range("N27").Value = -1210.6667
$ws.Range("H40").Value = 1889.4642
$ws.Range("I40").Value = 1456.5714
$ws.Range("K40").Value = 1456.5714
$ws.Range("M40").Value = -1320.5714
$ws.Range("H43").Value = 15666.667
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15666.667
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15666.667
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -16052.667
$ws.Range("H55").Value = 453.38235
$ws.Range("I55").Value = 442.375
$ws.Range("K55").Value = 442.375
$ws.Range("M55").Value = -269.375
$ws.Range("H93").Value = 13338368
$ws.Range("I93").Value = 16668497
$ws.Range("K93").Value = 16668497
$ws.Range("M93").Value = -16667249
$ws.Range("H100").Value = 66669344
$ws.Range("I100").Value = 111113550
$ws.Range("K100").Value = 111113550
$ws.Range("M100").Value = -111113009
$ws.Range("H122").Value = 7222.4287
$ws.Range("I122").Value = 6939.25
$ws.Range("J122").Value = 7600
$ws.Range("K122").Value = 20817.75
$ws.Range("L122").Value = 22800
$ws.Range("M122").Value = -18367.75
$ws.Range("N122").Value = -27700
$ws.Range("H132").Value = 40185.4
$ws.Range("I132").Value = 38467.656
$ws.Range("K132").Value = 115402.968
$ws.Range("M132").Value = -112872.968
$ws.Range("H136").Value = 7541.75
$ws.Range("I136").Value = 10367.5
$ws.Range("K136").Value = 31102.5
$ws.Range("M136").Value = -28552.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 626
$ws.Range("I100").Value = 422.3125
$ws.Range("K100").Value = 844.625
$ws.Range("M100").Value = -303.625
$ws.Range("H113").Value = 430.48572
$ws.Range("I113").Value = 425.61905
$ws.Range("K113").Value = 1276.85715
$ws.Range("M113").Value = 893.14285
$ws.Range("H122").Value = 8339.975
$ws.Range("I122").Value = 2785.4666
$ws.Range("J122").Value = 26855
$ws.Range("K122").Value = 8356.399800000001
$ws.Range("L122").Value = 80565
$ws.Range("M122").Value = -5906.399800000001
$ws.Range("N122").Value = -85465
